$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = 2
    3 = 1
    4 = 3
    5 = 0
    6 = 4
    7 = 6
    8 = 4
    9 = 1
    10 = 4
    11 = 0
    12 = 6
    13 = 5
    14 = 6
    15 = 5
    16 = 4
    17 = 1
    18 = 2
    19 = 2
    20 = 2
    21 = 2
    22 = 3
    23 = 5
    24 = 5
    25 = 4
    26 = 6
    27 = 6
    28 = 3
    29 = 11
    30 = 6
    31 = 3
    32 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
